# Project Backlog.xlsx - fill in remaining feature names for the backlog
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: clear the estimate inputs, keep the formula (now evaluates to 0) ---
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# --- Row 3 ---
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

# --- Row 4 ---
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# --- Row 6: add the missing Adjust Estimate formula ---
$ws.Range("E6").Formula = "=C6*(1+D6)"

# --- Row 8 ---
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

# --- Row 9: add the missing Adjust Estimate formula ---
$ws.Range("E9").Formula = "=C9*(1+D9)"

# --- Row 11: B11 was a placeholder "sdf"; rename it to the real feature name
#     that used to live in C11, then clear C11 and add the formula ---
$ws.Range("B11").Value = "Home Page Screen with Appropirate Buttons"
$ws.Range("C11").ClearContents()
$ws.Range("E11").Formula = "=C11*(1+D11)"

# --- Row 12: add the missing Adjust Estimate formula ---
$ws.Range("E12").Formula = "=C12*(1+D12)"

# --- Row 13: fill in the milestone id + formula ---
$ws.Range("A13").Value = "M6-3"
$ws.Range("E13").Formula = "=C13*(1+D13)"

# --- New rows for milestone M7 ---
$ws.Range("A15").Value = "M7-1"
$ws.Range("B15").Value = "Account Selection Screen"
$ws.Range("E15").Formula = "=C15*(1+D15)"

$ws.Range("B16").Value = "Account Balance Screen"
$ws.Range("E16").Formula = "=C16*(1+D16)"

$ws.Range("B17").Value = "New Transaction Screen with Appropriate Fields"
$ws.Range("E17").Formula = "=C17*(1+D17)"

$ws.Range("B18").Value = "New Transaction Handling"
$ws.Range("E18").Formula = "=C18*(1+D18)"

$ws.Range("A16").Value = "M7-2"
$ws.Range("A17").Value = "M7-3"
$ws.Range("A18").Value = "M7-4"

# --- New rows for milestone M8 ---
$ws.Range("A20").Value = "M8-1"
$ws.Range("E20").Formula = "=C20*(1+D20)"

$ws.Range("B21").Value = "Spending Category Report Screen"
$ws.Range("E21").Formula = "=C21*(1+D21)"

$ws.Range("B20").Value = "Spending Category Report Promt Screen"

$ws.Range("A21").Value = "M8-2"

# --- New rows for milestone M9 ---
$ws.Range("A23").Value = "M9-1"
$ws.Range("B23").Value = "Save/Loading Application"
$ws.Range("E23").Formula = "=C23*(1+D23)"

$ws.Range("A24").Value = "M9-2"
$ws.Range("B24").Value = "Local Information Storage"
$ws.Range("E24").Formula = "=C24*(1+D24)"

# --- New row for milestone M10 ---
$ws.Range("A26").Value = "M10"
$ws.Range("B26").Value = "Code Documentations"
$ws.Range("E26").Formula = "=C26*(1+D26)"

# --- Grand Total row ---
$ws.Range("A28").Value = "Total"
$ws.Range("C28").Formula = "=SUM(C2:C27)"
$ws.Range("E28").Formula = "=SUM(E2:E27)"
$ws.Range("G28").Formula = "=SUM(G2:G27)"
$ws.Range("H28").Formula = "=SUM(H2:H27)"
$ws.Range("I28").Formula = "=SUM(I2:I27)"

# --- Update the saved selection to match the author's final cursor position ---
$ws.Range("G12").Select()
